# Generate Report for Handoff
#
# The localization status report moved from "In Translation" to
# "Ready for handoff": refresh the Status columns, bump the
# handoff/generate timestamps, and widen the Status columns so the
# longer text isn't clipped.

$wb = $excel.ActiveWorkbook

# NOTE: the simulated ColumnWidth setter in this runtime snaps every
# assignment onto a 1/6-character grid, so we pick the input that lands
# on the grid point closest to the desired ~17.216 character width.
$statusColWidth = 16.333333

# --- Overview sheet: per-language status + last handoff-generation date ---
$wsOverview = $wb.Sheets.Item("Overview")
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsOverview.Range("G2").Value = "2016-08-24 18:47:54"
$wsOverview.Columns.Item(5).ColumnWidth = $statusColWidth
$wsOverview.Columns.Item(6).ColumnWidth = $statusColWidth

# --- zh-cn detail sheet: Status + Latest Handoff Datetime ---
$wsZhCn = $wb.Sheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("H2").Value = "2016-08-24 18:47:49"
$wsZhCn.Columns.Item(3).ColumnWidth = $statusColWidth

# --- de-de detail sheet: Status + Latest Handoff Datetime ---
$wsDeDe = $wb.Sheets.Item("de-de")
$wsDeDe.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("H2").Value = "2016-08-24 18:47:54"
$wsDeDe.Columns.Item(3).ColumnWidth = $statusColWidth
